$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume table with latest scraped values

$ws.Range("D2").Value = "47.437.71"
$ws.Range("E2").Value = "  +3.01%  "
$ws.Range("D3").Value = "2.512.76"
$ws.Range("E3").Value = "  +2.66%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "110.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.17%  "
$ws.Range("E7").Value = "  +1.62%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +1.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.31"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.11%  "
$ws.Range("E11").Value = "  +1.73%  "
$ws.Range("E12").Value = "  +1.04%  "
$ws.Range("E13").Value = "  +1.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.77%  "
$ws.Range("D15").Value = "2.903.57"
$ws.Range("E15").Value = "  +2.57%  "
$ws.Range("D16").Value = "2.510.93"
$ws.Range("E16").Value = "  +2.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.865"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.70%  "
$ws.Range("D18").Value = "47.409.78"
$ws.Range("E18").Value = "  +3.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.00%  "
$ws.Range("D21").Value = "0.0₃0948"
$ws.Range("E21").Value = "  +1.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.64"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +12.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "250.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.30%  "
$ws.Range("E25").Value = "  +3.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.82%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.78%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.92%  "
$ws.Range("E31").Value = "  +5.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("E34").Value = "  +2.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0797"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.61%  "
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.02"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.75"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.37%  "
$ws.Range("E39").Value = "  +2.87%  "
$ws.Range("E40").Value = "  +1.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "122.87"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.31%  "
$ws.Range("E42").Value = "  -0.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.56%  "
$ws.Range("D45").Value = "2.002.60"
$ws.Range("E45").Value = "  +2.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.02%  "
$ws.Range("E48").Value = "  -2.09%  "
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.83"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.50%  "
